# Update the "want to go" counts (column F) on the "展览" (Exhibition)
# sheet and the mirrored aggregate "全部类型" (All types) sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value  = 1071
$wsExhibit.Range("F10").Value = 128
$wsExhibit.Range("F12").Value = 154
$wsExhibit.Range("F14").Value = 2762
$wsExhibit.Range("F15").Value = 1037

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 1071
$wsAll.Range("F12").Value = 128
$wsAll.Range("F14").Value = 154
$wsAll.Range("F16").Value = 2762
$wsAll.Range("F17").Value = 1037
